$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.852.79"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "2.748.39"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'574.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").Value = "'157.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").Value = "'5.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -15.52%  "
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").Value = "3.235.00"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").Value = "'26.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").Value = "63.496.00"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "'0.0000150"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").Value = "2.753.47"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "'12.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("D20").Value = "'355.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("E21").Value = "  -3.87%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").Value = "'65.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").Value = "'0.170"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'8.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("D28").Value = "0.0₃0908"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("E29").Value = "  -4.00%  "
$ws.Range("D30").Value = "'7.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.23%  "
$ws.Range("D31").Value = "'1.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.08%  "
$ws.Range("D32").Value = "'168.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.53%  "
$ws.Range("D33").Value = "'20.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.16%  "
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").Value = "'1.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("D38").Value = "'0.979"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.07%  "
$ws.Range("E39").Value = "  +5.46%  "
$ws.Range("D40").Value = "'4.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("D41").Value = "'330.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("D42").Value = "'39.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").Value = "'21.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("D44").Value = "'0.0588"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").Value = "'21.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.07%  "
$ws.Range("D46").Value = "'0.0254"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'135.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.625"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.58%  "
$ws.Range("E49").Value = "  -1.19%  "
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").Value = "'11.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.09%  "
